# Generate Report for Handoff
# Adds a new tracked file (ac481a19-cbc4-42bf-abde-74a8bf5386cb.md) as a new
# row to the Overview / zh-cn / de-de tables.

$wb = $excel.ActiveWorkbook

$hyperlinkFontColor = 15570276  # BGR(ED,95,64) == RGB(0x64,0x95,0xED) == FF6495ED
$dateNumberFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkFontColor
}

function Style-AsDate($range) {
    $range.NumberFormat = $dateNumberFormat
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "ac481a19-cbc4-42bf-abde-74a8bf5386cb.md"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e6ab2b43cf3af615694cae559b00bc4a1091a06/e2e/ac481a19-cbc4-42bf-abde-74a8bf5386cb.md",
    "",
    "",
    "e2e\ac481a19-cbc4-42bf-abde-74a8bf5386cb.md"
) | Out-Null
Style-AsHyperlink $wsOverview.Range("B4")

$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2017-02-09 09:51:06"
Style-AsDate $wsOverview.Range("G4")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e6ab2b43cf3af615694cae559b00bc4a1091a06/e2e/ac481a19-cbc4-42bf-abde-74a8bf5386cb.md",
    "",
    "",
    "ac481a19-cbc4-42bf-abde-74a8bf5386cb.md"
) | Out-Null
Style-AsHyperlink $wsZhCn.Range("A4")

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "ac481a19-cbc4-42bf-abde-74a8bf5386cb.c7a89fd2c3843e97e05643dfcb1bb33563e00ba4.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2017-02-09 09:50:43"
Style-AsDate $wsZhCn.Range("H4")
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = ""
$wsZhCn.Range("L4").Value = "0001-01-01 00:00:00"
Style-AsDate $wsZhCn.Range("L4")
$wsZhCn.Range("M4").Value = ""
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "True"
$wsZhCn.Range("P4").Value = ""
$wsZhCn.Range("Q4").Value = "False"
$wsZhCn.Range("R4").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e6ab2b43cf3af615694cae559b00bc4a1091a06/e2e/ac481a19-cbc4-42bf-abde-74a8bf5386cb.md",
    "",
    "",
    "ac481a19-cbc4-42bf-abde-74a8bf5386cb.md"
) | Out-Null
Style-AsHyperlink $wsDeDe.Range("A4")

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "ac481a19-cbc4-42bf-abde-74a8bf5386cb.c7a89fd2c3843e97e05643dfcb1bb33563e00ba4.de-de.xlf"
$wsDeDe.Range("H4").Value = "2017-02-09 09:51:06"
Style-AsDate $wsDeDe.Range("H4")
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = ""
$wsDeDe.Range("L4").Value = "0001-01-01 00:00:00"
Style-AsDate $wsDeDe.Range("L4")
$wsDeDe.Range("M4").Value = ""
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "True"
$wsDeDe.Range("P4").Value = ""
$wsDeDe.Range("Q4").Value = "False"
$wsDeDe.Range("R4").Value = ""

Write-Host "Done adding handoff row for ac481a19-cbc4-42bf-abde-74a8bf5386cb.md"
